$wb = $excel.ActiveWorkbook

$groupWs = $wb.Worksheets.Item("Group")
$inputWs = $wb.Worksheets.Item("Input")

# Insert a new row for the "staff" / "Nail, Wilson" test step above the
# existing "GroupActive" row, shifting the remaining rows down.
$groupWs.Rows("3:3").Insert()
$groupWs.Range("A3").Value = "staff"
$groupWs.Range("B3").Value = "Nail, Wilson"

# Match the formatting used for the equivalent "staff" label/value pair
# elsewhere in the workbook (Input sheet).
$inputWs.Range("A11").Copy()
$groupWs.Range("A3").PasteSpecial(-4122)
$inputWs.Range("B1").Copy()
$groupWs.Range("B3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the saved selection on the Input sheet (no longer the active tab).
$inputWs.Range("B2").Select()

# Finally select the Group sheet, making it the active/visible tab with
# the cursor parked on B13.
$groupWs.Range("B13").Select()
